# Natmi following Dr Hou advice
# Replace the LR-pairs result table (rows 2-17, columns A:T) with the
# updated natmi output (re-computed with n_bins settings per Dr Hou's advice).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is one data row in column order A..T
$rows = @(
    @("ECs", "App", "Tnfrsf21", "ECs", 3, 1, 110.8604276666667, 332.581283, 0.2509786052589675, 0.2509786052589675, 3, 1, 15.44356366666667, 46.330691, 0.2921725357712574, 0.2921725357712575, 1712.080072784061, 15408.72065505655, 0.07332905552284597, 0.07332905552284598),
    @("ECs", "App", "Tnfrsf21", "FAPs", 3, 1, 110.8604276666667, 332.581283, 0.2509786052589675, 0.2509786052589675, 2, 0.6666666666666666, 1.215543, 3.646629, 0.02299652387975402, 0.02299652387975402, 134.755616827223, 1212.800551445007, 0.005771635489145204, 0.005771635489145204),
    @("ECs", "App", "Tnfrsf21", "M2", 3, 1, 110.8604276666667, 332.581283, 0.2509786052589675, 0.2509786052589675, 3, 1, 23.02356833333333, 69.070705, 0.4355765604134518, 0.4355765604134519, 2552.402631846057, 22971.62368661452, 0.1093203976160665, 0.1093203976160665),
    @("ECs", "App", "Tnfrsf21", "sCs", 3, 1, 110.8604276666667, 332.581283, 0.2509786052589675, 0.2509786052589675, 3, 1, 13.17500933333333, 39.525028, 0.2492543799355367, 0.2492543799355367, 1460.587169205658, 13145.28452285092, 0.06255751663090978, 0.06255751663090978),
    @("FAPs", "App", "Tnfrsf21", "ECs", 3, 1, 184.841802, 554.525406, 0.4184661617850055, 0.4184661617850055, 3, 1, 15.44356366666667, 46.330691, 0.2921725357712574, 0.2921725357712575, 2854.616137448394, 25691.54523703555, 0.1222643196231903, 0.1222643196231903),
    @("FAPs", "App", "Tnfrsf21", "FAPs", 3, 1, 184.841802, 554.525406, 0.4184661617850055, 0.4184661617850055, 2, 0.6666666666666666, 1.215543, 3.646629, 0.02299652387975402, 0.02299652387975402, 224.683158528486, 2022.148426756374, 0.009623267082357887, 0.009623267082357887),
    @("FAPs", "App", "Tnfrsf21", "M2", 3, 1, 184.841802, 554.525406, 0.4184661617850055, 0.4184661617850055, 3, 1, 23.02356833333333, 69.070705, 0.4355765604134518, 0.4355765604134519, 4255.71785920347, 38301.46073283123, 0.1822740513997317, 0.1822740513997318),
    @("FAPs", "App", "Tnfrsf21", "sCs", 3, 1, 184.841802, 554.525406, 0.4184661617850055, 0.4184661617850055, 3, 1, 13.17500933333333, 39.525028, 0.2492543799355367, 0.2492543799355367, 2435.292466540152, 21917.63219886137, 0.1043045236797255, 0.1043045236797255),
    @("M2", "App", "Tnfrsf21", "ECs", 3, 1, 95.23175666666667, 285.69527, 0.2155966197102082, 0.2155966197102082, 3, 1, 15.44356366666667, 46.330691, 0.2921725357712574, 0.2921725357712575, 1470.717697170175, 13236.45927453157, 0.06299141108444299, 0.06299141108444301),
    @("M2", "App", "Tnfrsf21", "FAPs", 3, 1, 95.23175666666667, 285.69527, 0.2155966197102082, 0.2155966197102082, 2, 0.6666666666666666, 1.215543, 3.646629, 0.02299652387975402, 0.02299652387975402, 115.75829519387, 1041.82465674483, 0.00495797281356005, 0.004957972813560051),
    @("M2", "App", "Tnfrsf21", "M2", 3, 1, 95.23175666666667, 285.69527, 0.2155966197102082, 0.2155966197102082, 3, 1, 23.02356833333333, 69.070705, 0.4355765604134518, 0.4355765604134519, 2192.574857118372, 19733.17371406535, 0.0939088340501395, 0.09390883405013953),
    @("M2", "App", "Tnfrsf21", "sCs", 3, 1, 95.23175666666667, 285.69527, 0.2155966197102082, 0.2155966197102082, 3, 1, 13.17500933333333, 39.525028, 0.2492543799355367, 0.2492543799355367, 1254.679282913062, 11292.11354621756, 0.05373840176206566, 0.05373840176206567),
    @("sCs", "App", "Tnfrsf21", "ECs", 3, 1, 50.778675, 152.336025, 0.1149586132458188, 0.1149586132458188, 3, 1, 15.44356366666667, 46.330691, 0.2921725357712574, 0.2921725357712575, 784.2037002714751, 7057.833302443276, 0.03358774954077813, 0.03358774954077814),
    @("sCs", "App", "Tnfrsf21", "FAPs", 3, 1, 50.778675, 152.336025, 0.1149586132458188, 0.1149586132458188, 2, 0.6666666666666666, 1.215543, 3.646629, 0.02299652387975402, 0.02299652387975402, 61.723662945525, 555.5129665097251, 0.002643648494690879, 0.002643648494690879),
    @("sCs", "App", "Tnfrsf21", "M2", 3, 1, 50.778675, 152.336025, 0.1149586132458188, 0.1149586132458188, 3, 1, 23.02356833333333, 69.070705, 0.4355765604134518, 0.4355765604134519, 1169.106293738625, 10521.95664364763, 0.05007327734751402, 0.05007327734751404),
    @("sCs", "App", "Tnfrsf21", "sCs", 3, 1, 50.778675, 152.336025, 0.1149586132458188, 0.1149586132458188, 3, 1, 13.17500933333333, 39.525028, 0.2492543799355367, 0.2492543799355367, 669.0095170593, 6021.0856535337, 0.02865393786283574, 0.02865393786283574)
)

$numRows = $rows.Count
$numCols = 20

$data = New-Object 'object[,]' $numRows,$numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $data[$r,$c] = $rows[$r][$c]
    }
}

$ws.Range("A2:T17").Value = $data
